# Insert a new "ProductNumber" column (with value "G3U45A") before the
# existing "QuoteName" column. This shifts the former columns F:J one
# position to the right (becoming G:K) and adds the new data in column F.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shift row 1 (header) cells one column to the right, right-to-left so
#     we never clobber a value before it has been copied onward. ---
$ws.Range("K1").Value = $ws.Range("J1").Value2
$ws.Range("J1").Value = $ws.Range("I1").Value2
$ws.Range("I1").Value = $ws.Range("H1").Value2
$ws.Range("H1").Value = $ws.Range("G1").Value2
$ws.Range("G1").Value = $ws.Range("F1").Value2
$ws.Range("F1").Value = "ProductNumber"

# --- Shift row 2 (data) cells one column to the right as well. ---
$ws.Range("J2").Value = $ws.Range("I2").Value2
$ws.Range("I2").Value = $ws.Range("H2").Value2
$ws.Range("H2").Value = $ws.Range("G2").Value2
$ws.Range("G2").Value = $ws.Range("F2").Value2
$ws.Range("F2").Value = "G3U45A"

# --- Shift the column widths to match, right-to-left, then give the new
#     ProductNumber column its own best-fit-like width. ---
$ws.Columns("K").ColumnWidth = 15.451822916666668
$ws.Columns("J").ColumnWidth = 16.736979166666668
$ws.Columns("I").ColumnWidth = 13.307291666666668
$ws.Columns("H").ColumnWidth = 22.592447916666668
$ws.Columns("G").ColumnWidth = 11.022135416666668
$ws.Columns("F").ColumnWidth = 14.451822916666668

# --- Update the active selection to reflect the new data entry cell. ---
$ws.Range("F2").Select() | Out-Null
